$d = $word.ActiveDocument

# Date header
$d.Content.Find.Execute("2025-03-16 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-17 Monday", 2)

# Table cell values (all unique strings in the document).
# NOTE: "788×8=6304" is both a pre-existing value (originally paired with
# "401×4=1604" row) and the *new* value that "806×9=7254" becomes. To avoid
# Find/Replace-All clobbering both occurrences, the original "788×8=6304"
# must be replaced BEFORE the "806×9=7254" -> "788×8=6304" replacement runs.

$d.Content.Find.Execute("553×5=2765", $true, $false, $false, $false, $false, $true, 1, $false, "891×6=5346", 2)
$d.Content.Find.Execute("568×4=2272", $true, $false, $false, $false, $false, $true, 1, $false, "309×8=2472", 2)
$d.Content.Find.Execute("404×6=2424", $true, $false, $false, $false, $false, $true, 1, $false, "639×4=2556", 2)
$d.Content.Find.Execute("650×8=5200", $true, $false, $false, $false, $false, $true, 1, $false, "166×2=332", 2)
$d.Content.Find.Execute("151×4=604", $true, $false, $false, $false, $false, $true, 1, $false, "698×5=3490", 2)

$d.Content.Find.Execute("811×7=5677", $true, $false, $false, $false, $false, $true, 1, $false, "960×6=5760", 2)
$d.Content.Find.Execute("167×7=1169", $true, $false, $false, $false, $false, $true, 1, $false, "359×2=718", 2)
$d.Content.Find.Execute("488×5=2440", $true, $false, $false, $false, $false, $true, 1, $false, "215×7=1505", 2)
$d.Content.Find.Execute("810×7=5670", $true, $false, $false, $false, $false, $true, 1, $false, "923×5=4615", 2)
$d.Content.Find.Execute("603×2=1206", $true, $false, $false, $false, $false, $true, 1, $false, "937×8=7496", 2)

$d.Content.Find.Execute("900×4=3600", $true, $false, $false, $false, $false, $true, 1, $false, "714×7=4998", 2)
$d.Content.Find.Execute("178×2=356", $true, $false, $false, $false, $false, $true, 1, $false, "368×8=2944", 2)
$d.Content.Find.Execute("997×3=2991", $true, $false, $false, $false, $false, $true, 1, $false, "786×3=2358", 2)
$d.Content.Find.Execute("222×3=666", $true, $false, $false, $false, $false, $true, 1, $false, "812×9=7308", 2)
$d.Content.Find.Execute("937×2=1874", $true, $false, $false, $false, $false, $true, 1, $false, "257×4=1028", 2)

$d.Content.Find.Execute("697×8=5576", $true, $false, $false, $false, $false, $true, 1, $false, "419×9=3771", 2)
$d.Content.Find.Execute("367×5=1835", $true, $false, $false, $false, $false, $true, 1, $false, "149×4=596", 2)
$d.Content.Find.Execute("788×8=6304", $true, $false, $false, $false, $false, $true, 1, $false, "855×7=5985", 2)
$d.Content.Find.Execute("401×4=1604", $true, $false, $false, $false, $false, $true, 1, $false, "631×4=2524", 2)
$d.Content.Find.Execute("187×8=1496", $true, $false, $false, $false, $false, $true, 1, $false, "812×6=4872", 2)

$d.Content.Find.Execute("263×4=1052", $true, $false, $false, $false, $false, $true, 1, $false, "305×3=915", 2)
$d.Content.Find.Execute("285×7=1995", $true, $false, $false, $false, $false, $true, 1, $false, "929×8=7432", 2)
$d.Content.Find.Execute("504×2=1008", $true, $false, $false, $false, $false, $true, 1, $false, "189×6=1134", 2)
$d.Content.Find.Execute("232×7=1624", $true, $false, $false, $false, $false, $true, 1, $false, "723×7=5061", 2)
$d.Content.Find.Execute("806×9=7254", $true, $false, $false, $false, $false, $true, 1, $false, "788×8=6304", 2)
